$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3,3,3,1,1,2,3,3,2,1,2,3,3,1,3,2,3,3,3,3,1,1,2,3,3,3,2,1,1,3,3,3,3,3,3,3,3,3,1,2,3,3,3,3,1,3,3,3,1,1,3,1,2,1,1,3,1,2,2,3,3,2,1,3,1,2,3,3,1,1,1,3,2,3,3,2,1,1,3,1,1,3,1,2,3,2,3,2,1,3,3,1,2,2,1,2,1,3,1,3,1,3,3,1,3,3,3,3,1,1,2,1,1,2,3,1,3,1,3,1,3,1,1,1,3,3,1,1,2,1,3,3,3,3,1,3,3,3,1,1,2,3,2,3,3,3,3,2,3,1,3,1,1,1,2,1,3,1,2,3,2,2,1,1,1,1,3,2,1,1,3,1,1,3,1,1,1,2,1,1,3,3,3,3,3,3,3,2,1,3,3,3,3,1,3,1,2,2,3,3,1,2,2,3,1,1,1,2,3,3,2,1,1,3,3,3,3,3,2,3,1,3,1,1,1,1,3,2,1,3,1,3,1,3,3,2,3,1,1,3,3,3,2,3,3,3,1,3,2,2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = [int]$values[$i]
}
